$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 15: new scaffolding test entry for day meter card
$ws.Range("A15").Value = "high"
$ws.Range("B15").Value = "Day meter card contains day stats"

# Row 16: new scaffolding test entry for calendar directive wrapping
$ws.Range("A16").Value = "high"
$ws.Range("B16").Value = "Calendar is wrapped in a directive"

# Update selection to match the post-edit state (A16 selected)
$ws.Range("A16").Select()
